$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "positive integer"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = "PASS"

# New row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "negative"
$ws.Range("C3").Value = -9
$ws.Range("D3").Value = -204
$ws.Range("E3").Value = -204
$ws.Range("F3").Value = "PASS"

# New row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "help command"
$ws.Range("C4").Value = "argv[0] -h"
$ws.Range("D4").Value = "Input:Integer"
$ws.Range("E4").Value = "Input:Integer"
$ws.Range("F4").Value = "PASS"

# Update selection to F4 like the diff shows
$ws.Range("F4").Select()
